# Scheduled runner update: refresh market-board price/profit figures for
# several Leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 10870728
$ws.Range("I15").Value = 10870728
$ws.Range("K15").Value = 32612184
$ws.Range("M15").Value = -32612015

$ws.Range("H18").Value = 55558850
$ws.Range("I18").Value = 66668850
$ws.Range("J18").Value = 8850
$ws.Range("K18").Value = 66668850
$ws.Range("L18").Value = 8850
$ws.Range("M18").Value = -66668566
$ws.Range("N18").Value = -9418

$ws.Range("H116").Value = 34642064
$ws.Range("I116").Value = 50702908
$ws.Range("K116").Value = 50702908
$ws.Range("M116").Value = -50699466

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5353.392
$ws.Range("J32").Value = 16043.3
$ws.Range("L32").Value = 16043.3
$ws.Range("N32").Value = -16617.3

$ws.Range("H45").Value = 2882.4348
$ws.Range("I45").Value = 2621.1052
$ws.Range("J45").Value = 4123.75
$ws.Range("K45").Value = 2621.1052
$ws.Range("L45").Value = 4123.75
$ws.Range("M45").Value = -2244.1052
$ws.Range("N45").Value = -4877.75

$ws.Range("H61").Value = 17479.4
$ws.Range("I61").Value = 26883.166
$ws.Range("K61").Value = 26883.166
$ws.Range("M61").Value = -26671.166

$ws.Range("H74").Value = 8930156
$ws.Range("I74").Value = 12501125
$ws.Range("K74").Value = 12501125
$ws.Range("M74").Value = -12500251

$ws.Range("H77").Value = 8930156
$ws.Range("I77").Value = 12501125
$ws.Range("K77").Value = 62505625
$ws.Range("M77").Value = -62501257

$ws.Range("H97").Value = 906.625
$ws.Range("I97").Value = 1125.5
$ws.Range("K97").Value = 1125.5
$ws.Range("M97").Value = -629.5

$ws.Range("H122").Value = 3080.5
$ws.Range("I122").Value = 1604.7
$ws.Range("K122").Value = 4814.1
$ws.Range("M122").Value = -2364.1

$ws.Range("H132").Value = 5816.1333
$ws.Range("I132").Value = 2278.1052
$ws.Range("K132").Value = 6834.3156
$ws.Range("M132").Value = -4304.3156

$ws.Range("H136").Value = 17479.4
$ws.Range("I136").Value = 26883.166
$ws.Range("K136").Value = 80649.49800000001
$ws.Range("M136").Value = -78099.49800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2229.875
$ws.Range("I86").Value = 2183.1667
$ws.Range("J86").Value = 2370
$ws.Range("K86").Value = 2183.1667
$ws.Range("L86").Value = 2370
$ws.Range("M86").Value = -1060.1667
$ws.Range("N86").Value = -4616

$ws.Range("H89").Value = 2229.875
$ws.Range("I89").Value = 2183.1667
$ws.Range("J89").Value = 2370
$ws.Range("K89").Value = 10915.8335
$ws.Range("L89").Value = 11850
$ws.Range("M89").Value = -5299.833500000001
$ws.Range("N89").Value = -23082

$ws.Range("H107").Value = 1599.7142
$ws.Range("I107").Value = 1581.2222
$ws.Range("J107").Value = 1710.6666
$ws.Range("K107").Value = 1581.2222
$ws.Range("L107").Value = 1710.6666
$ws.Range("M107").Value = 338.7778000000001
$ws.Range("N107").Value = -5550.6666

$ws.Range("H134").Value = 6443.1113
$ws.Range("I134").Value = 1999.75
$ws.Range("K134").Value = 5999.25
$ws.Range("M134").Value = -3464.25

$ws.Range("H137").Value = 40000
$ws.Range("J137").Value = 40000
$ws.Range("L137").Value = 40000
$ws.Range("N137").Value = -50200

$ws.Range("H140").Value = 230000
$ws.Range("J140").Value = 230000
$ws.Range("L140").Value = 230000
$ws.Range("N140").Value = -240360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5915.871
$ws.Range("I31").Value = 2171.7856
$ws.Range("K31").Value = 2171.7856
$ws.Range("M31").Value = -1876.7856

$ws.Range("H34").Value = 5915.871
$ws.Range("I34").Value = 2171.7856
$ws.Range("K34").Value = 2171.7856
$ws.Range("M34").Value = -1969.7856

$ws.Range("H58").Value = 629083.5
$ws.Range("I58").Value = 718288.4399999999
$ws.Range("J58").Value = 4649
$ws.Range("K58").Value = 718288.4399999999
$ws.Range("L58").Value = 4649
$ws.Range("M58").Value = -718085.4399999999
$ws.Range("N58").Value = -5055

$ws.Range("H107").Value = 1515637.1
$ws.Range("I107").Value = 1818589.6
$ws.Range("K107").Value = 1818589.6
$ws.Range("M107").Value = -1816669.6

$ws.Range("H132").Value = 16684924
$ws.Range("I132").Value = 20848906
$ws.Range("J132").Value = 28999.25
$ws.Range("K132").Value = 62546718
$ws.Range("L132").Value = 86997.75
$ws.Range("M132").Value = -62544188
$ws.Range("N132").Value = -92057.75

$ws.Range("H134").Value = 4088.158
$ws.Range("I134").Value = 4104.8125
$ws.Range("K134").Value = 12314.4375
$ws.Range("M134").Value = -9779.4375

$ws.Range("H136").Value = 629083.5
$ws.Range("I136").Value = 718288.4399999999
$ws.Range("J136").Value = 4649
$ws.Range("K136").Value = 2154865.32
$ws.Range("L136").Value = 13947
$ws.Range("M136").Value = -2152315.32
$ws.Range("N136").Value = -19047

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 224.2
$ws.Range("J12").Value = 229.71428
$ws.Range("L12").Value = 689.14284
$ws.Range("N12").Value = -1035.14284

$ws.Range("H17").Value = 1050
$ws.Range("I17").Value = 1050
$ws.Range("K17").Value = 3150
$ws.Range("M17").Value = -2981

$ws.Range("H34").Value = 3437093.2
$ws.Range("I34").Value = 5154640
$ws.Range("J34").Value = 2000
$ws.Range("K34").Value = 15463920
$ws.Range("L34").Value = 6000
$ws.Range("M34").Value = -15463836
$ws.Range("N34").Value = -6168

$ws.Range("H61").Value = 555.5
$ws.Range("I61").Value = 374
$ws.Range("J61").Value = 633.2857
$ws.Range("K61").Value = 1122
$ws.Range("L61").Value = 1899.8571
$ws.Range("M61").Value = -907
$ws.Range("N61").Value = -2329.8571

$ws.Range("H114").Value = 3304.0908
$ws.Range("I114").Value = 669
$ws.Range("K114").Value = 2007
$ws.Range("M114").Value = 1247

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 677910.2
$ws.Range("I80").Value = 1666976.9
$ws.Range("K80").Value = 1666976.9
$ws.Range("M80").Value = -1665978.9

$ws.Range("H83").Value = 677910.2
$ws.Range("I83").Value = 1666976.9
$ws.Range("K83").Value = 8334884.5
$ws.Range("M83").Value = -8329892.5

$ws.Range("H132").Value = 9198.143
$ws.Range("I132").Value = 9064.5
$ws.Range("K132").Value = 27193.5
$ws.Range("M132").Value = -24663.5

$ws.Range("H134").Value = 916784
$ws.Range("J134").Value = 916784
$ws.Range("L134").Value = 2750352
$ws.Range("N134").Value = -2755422

$ws.Range("H136").Value = 73275.336
$ws.Range("J136").Value = 73275.336
$ws.Range("L136").Value = 219826.008
$ws.Range("N136").Value = -224926.008

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 8002275.5
$ws.Range("I16").Value = 9092869
$ws.Range("K16").Value = 9092869
$ws.Range("M16").Value = -9092699

$ws.Range("H122").Value = 39412950
$ws.Range("J122").Value = 11909690
$ws.Range("L122").Value = 35729070
$ws.Range("N122").Value = -35733970

$ws.Range("H135").Value = 81764
$ws.Range("J135").Value = 81764
$ws.Range("L135").Value = 81764
$ws.Range("N135").Value = -91904

$ws.Range("H136").Value = 2353.9333
$ws.Range("I136").Value = 2118.4443
$ws.Range("K136").Value = 6355.3329
$ws.Range("M136").Value = -3805.3329

$ws.Range("H140").Value = 61189.35
$ws.Range("J140").Value = 61189.35
$ws.Range("L140").Value = 61189.35
$ws.Range("N140").Value = -71549.35000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2056.7144
$ws.Range("I107").Value = 2387.0454
$ws.Range("J107").Value = 845.5
$ws.Range("K107").Value = 7161.1362
$ws.Range("L107").Value = 2536.5
$ws.Range("M107").Value = -5241.1362
$ws.Range("N107").Value = -6376.5

$ws.Range("H132").Value = 111114060
$ws.Range("I132").Value = 55555556
$ws.Range("K132").Value = 166666668
$ws.Range("M132").Value = -166664138

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
